$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AugustRaw")

# Header row
$ws.Cells.Item(1, 1).Value = "Library"
$ws.Cells.Item(1, 2).Value = "Items owned by this library checked out at this library this month"
$ws.Cells.Item(1, 3).Value = "Items owned by other libraries checked out at this library this month"
$ws.Cells.Item(1, 4).Value = "Total circulation this month"

# Data rows
$ws.Cells.Item(2, 1).Value = "Atchison Public Library"
$ws.Cells.Item(2, 2).Value = 4609
$ws.Cells.Item(2, 3).Value = 1923
$ws.Cells.Item(2, 4).Value = 6532

$ws.Cells.Item(3, 1).Value = "Baldwin City Public Library"
$ws.Cells.Item(3, 2).Value = 3086
$ws.Cells.Item(3, 3).Value = 647
$ws.Cells.Item(3, 4).Value = 3733

$ws.Cells.Item(4, 1).Value = "Basehor Community Library"
$ws.Cells.Item(4, 2).Value = 10152
$ws.Cells.Item(4, 3).Value = 1397
$ws.Cells.Item(4, 4).Value = 11549

$ws.Cells.Item(5, 1).Value = "Bern Community Library"
$ws.Cells.Item(5, 2).Value = 126
$ws.Cells.Item(5, 3).Value = 73
$ws.Cells.Item(5, 4).Value = 199

$ws.Cells.Item(6, 1).Value = "Bonner Springs City Library"
$ws.Cells.Item(6, 2).Value = 5651
$ws.Cells.Item(6, 3).Value = 1247
$ws.Cells.Item(6, 4).Value = 6898

$ws.Cells.Item(7, 1).Value = "Burlingame Community Library"
$ws.Cells.Item(7, 2).Value = 442
$ws.Cells.Item(7, 3).Value = 148
$ws.Cells.Item(7, 4).Value = 590

$ws.Cells.Item(8, 1).Value = "Carbondale City Library"
$ws.Cells.Item(8, 2).Value = 701
$ws.Cells.Item(8, 3).Value = 148
$ws.Cells.Item(8, 4).Value = 849

$ws.Cells.Item(9, 1).Value = "Centralia Community Library"
$ws.Cells.Item(9, 2).Value = 322
$ws.Cells.Item(9, 3).Value = 33
$ws.Cells.Item(9, 4).Value = 355

$ws.Cells.Item(10, 1).Value = "Corning City Library"
$ws.Cells.Item(10, 2).Value = 46
$ws.Cells.Item(10, 4).Value = 46

$ws.Cells.Item(11, 1).Value = "Digital Content"

$ws.Cells.Item(12, 1).Value = "Doniphan County Library - Elwood"
$ws.Cells.Item(12, 2).Value = 182
$ws.Cells.Item(12, 3).Value = 12
$ws.Cells.Item(12, 4).Value = 194

$ws.Cells.Item(13, 1).Value = "Doniphan County Library - Highland"
$ws.Cells.Item(13, 2).Value = 259
$ws.Cells.Item(13, 3).Value = 201
$ws.Cells.Item(13, 4).Value = 460

$ws.Cells.Item(14, 1).Value = "Doniphan County Library - Troy"
$ws.Cells.Item(14, 2).Value = 560
$ws.Cells.Item(14, 3).Value = 142
$ws.Cells.Item(14, 4).Value = 702

$ws.Cells.Item(15, 1).Value = "Doniphan County Library - Wathena"
$ws.Cells.Item(15, 2).Value = 500
$ws.Cells.Item(15, 3).Value = 105
$ws.Cells.Item(15, 4).Value = 605

$ws.Cells.Item(16, 1).Value = "Effingham Community Library"
$ws.Cells.Item(16, 2).Value = 302
$ws.Cells.Item(16, 3).Value = 55
$ws.Cells.Item(16, 4).Value = 357

$ws.Cells.Item(17, 1).Value = "Eudora Community Library"
$ws.Cells.Item(17, 2).Value = 1820
$ws.Cells.Item(17, 3).Value = 650
$ws.Cells.Item(17, 4).Value = 2470

$ws.Cells.Item(18, 1).Value = "Everest, Barnes Reading Room"
$ws.Cells.Item(18, 2).Value = 127
$ws.Cells.Item(18, 3).Value = 22
$ws.Cells.Item(18, 4).Value = 149

$ws.Cells.Item(19, 1).Value = "Hiawatha, Morrill Public Library"
$ws.Cells.Item(19, 2).Value = 2166
$ws.Cells.Item(19, 3).Value = 638
$ws.Cells.Item(19, 4).Value = 2804

$ws.Cells.Item(20, 1).Value = "Highland Community College"
$ws.Cells.Item(20, 2).Value = 57
$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(20, 4).Value = 60

$ws.Cells.Item(21, 1).Value = "Holton, Beck-Bookman Library"
$ws.Cells.Item(21, 2).Value = 2054
$ws.Cells.Item(21, 3).Value = 516
$ws.Cells.Item(21, 4).Value = 2570

$ws.Cells.Item(22, 1).Value = "Horton Public Library"
$ws.Cells.Item(22, 2).Value = 110
$ws.Cells.Item(22, 3).Value = 34
$ws.Cells.Item(22, 4).Value = 144

$ws.Cells.Item(23, 1).Value = "Lansing Community Library"
$ws.Cells.Item(23, 2).Value = 2247
$ws.Cells.Item(23, 3).Value = 726
$ws.Cells.Item(23, 4).Value = 2973

$ws.Cells.Item(24, 1).Value = "Leavenworth Public Library"
$ws.Cells.Item(24, 2).Value = 9466
$ws.Cells.Item(24, 3).Value = 1822
$ws.Cells.Item(24, 4).Value = 11288

$ws.Cells.Item(25, 1).Value = "Linwood Community Library"
$ws.Cells.Item(25, 2).Value = 693
$ws.Cells.Item(25, 3).Value = 128
$ws.Cells.Item(25, 4).Value = 821

$ws.Cells.Item(26, 1).Value = "Louisburg Library"

$ws.Cells.Item(27, 1).Value = "Lyndon Carnegie Library"
$ws.Cells.Item(27, 2).Value = 509
$ws.Cells.Item(27, 3).Value = 300
$ws.Cells.Item(27, 4).Value = 809

$ws.Cells.Item(28, 1).Value = "McLouth Public Library"
$ws.Cells.Item(28, 2).Value = 139
$ws.Cells.Item(28, 3).Value = 66
$ws.Cells.Item(28, 4).Value = 205

$ws.Cells.Item(29, 1).Value = "Meriden-Ozawkie Public Library"
$ws.Cells.Item(29, 2).Value = 1522
$ws.Cells.Item(29, 3).Value = 539
$ws.Cells.Item(29, 4).Value = 2061

$ws.Cells.Item(30, 1).Value = "Northeast Kansas Library System"
$ws.Cells.Item(30, 2).Value = 16
$ws.Cells.Item(30, 3).Value = 55
$ws.Cells.Item(30, 4).Value = 71

$ws.Cells.Item(31, 1).Value = "Nortonville Public Library"
$ws.Cells.Item(31, 2).Value = 318
$ws.Cells.Item(31, 3).Value = 68
$ws.Cells.Item(31, 4).Value = 386

$ws.Cells.Item(32, 1).Value = "Osage City Library"
$ws.Cells.Item(32, 2).Value = 1824
$ws.Cells.Item(32, 3).Value = 539
$ws.Cells.Item(32, 4).Value = 2363

$ws.Cells.Item(33, 1).Value = "Osawatomie Public Library"
$ws.Cells.Item(33, 2).Value = 1119
$ws.Cells.Item(33, 3).Value = 456
$ws.Cells.Item(33, 4).Value = 1575

$ws.Cells.Item(34, 1).Value = "Oskaloosa Public Library"
$ws.Cells.Item(34, 2).Value = 556
$ws.Cells.Item(34, 3).Value = 226
$ws.Cells.Item(34, 4).Value = 782

$ws.Cells.Item(35, 1).Value = "Ottawa Library"
$ws.Cells.Item(35, 2).Value = 7526
$ws.Cells.Item(35, 3).Value = 1238
$ws.Cells.Item(35, 4).Value = 8764

$ws.Cells.Item(36, 1).Value = "Overbrook Public Library"
$ws.Cells.Item(36, 2).Value = 942
$ws.Cells.Item(36, 3).Value = 209
$ws.Cells.Item(36, 4).Value = 1151

$ws.Cells.Item(37, 1).Value = "Paola Free Library"
$ws.Cells.Item(37, 2).Value = 3113
$ws.Cells.Item(37, 3).Value = 502
$ws.Cells.Item(37, 4).Value = 3615

$ws.Cells.Item(38, 1).Value = "Perry-Lecompton Community Library"
$ws.Cells.Item(38, 2).Value = 121
$ws.Cells.Item(38, 3).Value = 13
$ws.Cells.Item(38, 4).Value = 134

$ws.Cells.Item(39, 1).Value = "Pomona Community Library"
$ws.Cells.Item(39, 2).Value = 153
$ws.Cells.Item(39, 3).Value = 70
$ws.Cells.Item(39, 4).Value = 223

$ws.Cells.Item(40, 1).Value = "Prairie Hills Schools - Axtell Public School"
$ws.Cells.Item(40, 2).Value = 497
$ws.Cells.Item(40, 3).Value = 9
$ws.Cells.Item(40, 4).Value = 506

$ws.Cells.Item(41, 1).Value = "Prairie Hills Schools - Sabetha Elementary School"
$ws.Cells.Item(41, 2).Value = 1448
$ws.Cells.Item(41, 3).Value = 24
$ws.Cells.Item(41, 4).Value = 1472

$ws.Cells.Item(42, 1).Value = "Prairie Hills Schools - Sabetha High School"
$ws.Cells.Item(42, 2).Value = 62
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 4).Value = 63

$ws.Cells.Item(43, 1).Value = "Prairie Hills Schools - Sabetha Middle School"
$ws.Cells.Item(43, 2).Value = 255
$ws.Cells.Item(43, 3).Value = 19
$ws.Cells.Item(43, 4).Value = 274

$ws.Cells.Item(44, 1).Value = "Prairie Hills Schools - Wetmore Academic Center (Permanently closed)"

$ws.Cells.Item(45, 1).Value = "Richmond Public Library"
$ws.Cells.Item(45, 2).Value = 257
$ws.Cells.Item(45, 3).Value = 70
$ws.Cells.Item(45, 4).Value = 327

$ws.Cells.Item(46, 1).Value = "Rossville Community Library"
$ws.Cells.Item(46, 2).Value = 1509
$ws.Cells.Item(46, 3).Value = 543
$ws.Cells.Item(46, 4).Value = 2052

$ws.Cells.Item(47, 1).Value = "Sabetha, Mary Cotton Library"
$ws.Cells.Item(47, 2).Value = 2996
$ws.Cells.Item(47, 3).Value = 859
$ws.Cells.Item(47, 4).Value = 3855

$ws.Cells.Item(48, 1).Value = "Seneca Free Library"
$ws.Cells.Item(48, 2).Value = 1816
$ws.Cells.Item(48, 3).Value = 227
$ws.Cells.Item(48, 4).Value = 2043

$ws.Cells.Item(49, 1).Value = "Silver Lake Library"
$ws.Cells.Item(49, 2).Value = 1376
$ws.Cells.Item(49, 3).Value = 394
$ws.Cells.Item(49, 4).Value = 1770

$ws.Cells.Item(50, 1).Value = "Tonganoxie Public Library"
$ws.Cells.Item(50, 2).Value = 3624
$ws.Cells.Item(50, 3).Value = 991
$ws.Cells.Item(50, 4).Value = 4615

$ws.Cells.Item(51, 1).Value = "Valley Falls, Delaware Township Library"
$ws.Cells.Item(51, 2).Value = 505
$ws.Cells.Item(51, 3).Value = 246
$ws.Cells.Item(51, 4).Value = 751

$ws.Cells.Item(52, 1).Value = "Wellsville City Library"
$ws.Cells.Item(52, 2).Value = 1145
$ws.Cells.Item(52, 3).Value = 347
$ws.Cells.Item(52, 4).Value = 1492

$ws.Cells.Item(53, 1).Value = "Wetmore Public Library"
$ws.Cells.Item(53, 2).Value = 155
$ws.Cells.Item(53, 3).Value = 154
$ws.Cells.Item(53, 4).Value = 309

$ws.Cells.Item(54, 1).Value = "Williamsburg Community Library"
$ws.Cells.Item(54, 2).Value = 266
$ws.Cells.Item(54, 3).Value = 32
$ws.Cells.Item(54, 4).Value = 298

$ws.Cells.Item(55, 1).Value = "Winchester Public Library"
$ws.Cells.Item(55, 2).Value = 351
$ws.Cells.Item(55, 3).Value = 259
$ws.Cells.Item(55, 4).Value = 610

